$d = $word.ActiveDocument
$sec = $d.Sections(1)
$footer = $sec.Footers(1)
$footer.Range.Tables.Add($footer.Range, 1, 3)
